$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Tarifa de equilíbrio"

# Header
$ws.Range("B1").Value = "Valores"

# Row labels
$ws.Range("A2").Value = "T-vig"
$ws.Range("A3").Value = "T-eq"

# Values (updated figures) with currency number format
$ws.Range("B2").Value = 4.05
$ws.Range("B3").Value = 4.29845246718807

$ws.Range("B2:B3").NumberFormat = "R$ #.##0,0"

# The labels no longer carry the old thin-box border
$ws.Range("A2").Borders.LineStyle = 0
$ws.Range("A3").Borders.LineStyle = 0
$ws.Range("B1").Borders.LineStyle = 0
